$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")
$ws.Range("E1").Comment.Delete()
$ws.Range("E1").EntireColumn.Delete()
[void]$ws.Range("F9").Select()
